# ValueSet-cbc-blood-automated-vs.xlsx : refresh IG publication metadata.
#   - Version 0.1.6 -> 0.1.7
#   - Status  active -> draft
#   - Date    2023-05-05T10:50:04-05:00 -> 2024-08-27T12:23:18-05:00
#   - Contact: publisher org now includes its URL, and a second "Contact"
#     row (the responsible author) plus a new "Jurisdiction" row are added.
# The "Include from LOINC" sheet (the LOINC concept list) is unchanged.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")

# --- Make room for two new rows by shifting rows 12-15 down to 13-16.
#     Copy bottom-up so each source row is read before it gets overwritten;
#     using Range.Copy (rather than setting .Value) carries the existing
#     border/fill/wrap style ("s=2") along instead of minting a new style.
$ws1.Range("A15:B15").Copy($ws1.Range("A16:B16"))
$ws1.Range("A14:B14").Copy($ws1.Range("A15:B15"))
$ws1.Range("A13:B13").Copy($ws1.Range("A14:B14"))
$ws1.Range("A12:B12").Copy($ws1.Range("A13:B13"))

# Purpose/Copyright (now rows 14 & 15) have no Value column - make sure no
# stale copied text lingers there.
$ws1.Range("B14").ClearContents()
$ws1.Range("B15").ClearContents()

# --- In-place metadata updates ---
$ws1.Range("B3").Value  = "0.1.7"                                    # Version
$ws1.Range("B6").Value  = "draft"                                    # Status
$ws1.Range("B8").Value  = "2024-08-27T12:23:18-05:00"                # Date
$ws1.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"  # Contact (publisher)

# --- New row 11: second Contact entry (the IG's responsible author) ---
$ws1.Range("A11").Value = "Contact"
$ws1.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# --- New row 12: Jurisdiction property, left without a value ---
$ws1.Range("A12").Value = "Jurisdiction"
$ws1.Range("B12").ClearContents()

$ws1.Range("A1").Select()
